$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# The "Rule" label in B11 (row for R40) changes from the text "R40" to the
# text "1". We need the result to stay a genuine text cell (it was typed as
# text before too), not get auto-coerced into a number, and we must not
# disturb B11's existing cell style (s="23").
#
# Typing an apostrophe-prefixed numeric string directly into B11 would force
# Excel to keep it as text, but it also stamps the cell with a brand new
# "quote prefixed" style variant. To avoid that we stage the text on a
# scratch cell that already carries the exact same style as B11 used to
# have, copy only its *value* onto B11 (so B11's original style survives
# untouched), and then restore the scratch cell's original formatting
# (borrowed from another cell sharing that same style) so nothing else in
# the sheet is left changed.
$scratch = $ws.Cells.Item(5, 2)    # B5 - currently empty, style matches B11's border style family
$formatDonor = $ws.Cells.Item(5, 5) # E5 - same original style as B5, used to restore it afterwards

$scratch.Value = "'1"
$scratch.Copy()

$target = $ws.Cells.Item(11, 2)    # B11
$target.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)

$scratch.ClearContents()
$formatDonor.Copy()
$scratch.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
